$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('K2').Value = 5048
$ws.Range('K3').Value = 5172
$ws.Range('K4').Value = 1074
$ws.Range('K5').Value = 368
$ws.Range('K6').Value = 5808
$ws.Range('K7').Value = 17470
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('K7').Value = 522
$ws.Range('K8').Value = 1170
$ws.Range('K10').Value = 95
$ws.Range('K11').Value = 338
$ws.Range('K13').Value = 19
$ws.Range('K15').Value = 176
$ws.Range('K18').Value = 118
$ws.Range('K20').Value = 399
$ws.Range('K23').Value = 176
$ws.Range('K25').Value = 83
$ws.Range('K27').Value = 161
$ws.Range('K29').Value = 940
$ws.Range('K31').Value = 191
$ws.Range('K33').Value = 741
$ws.Range('K34').Value = 94
$ws.Range('K36').Value = 230
$ws.Range('K37').Value = 589
$ws.Range('K42').Value = 647
$ws.Range('K44').Value = 155
$ws.Range('K49').Value = 99
$ws.Range('K52').Value = 454
$ws.Range('K54').Value = 342
$ws.Range('K55').Value = 199
$ws.Range('K59').Value = 29
$ws.Range('K63').Value = 51
$ws.Range('K65').Value = 404
$ws.Range('K67').Value = 671
$ws.Range('K73').Value = 151
$ws.Range('K76').Value = 241
$ws.Range('K77').Value = 125
$ws.Range('K78').Value = 201
$ws.Range('K80').Value = 63
$ws.Range('K81').Value = 12
$ws.Range('K85').Value = 810
$ws.Range('K88').Value = 194
$ws.Range('K89').Value = 254
$ws.Range('K90').Value = 155
$ws.Range('K91').Value = 187
$ws.Range('K94').Value = 231
$ws.Range('K96').Value = 188
$ws.Range('K101').Value = 17470
$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('K2').Value = 60
$ws.Range('K3').Value = 35
$ws.Range('K7').Value = 188
$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('K2').Value = 180
$ws.Range('K3').Value = 170
$ws.Range('K6').Value = 134
$ws.Range('K7').Value = 522
$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('K2').Value = 112
$ws.Range('K3').Value = 89
$ws.Range('K7').Value = 338
$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('K2').Value = 70
$ws.Range('K4').Value = 29
$ws.Range('K6').Value = 77
$ws.Range('K7').Value = 254
$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('K2').Value = 270
$ws.Range('K3').Value = 273
$ws.Range('K7').Value = 810
$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('K3').Value = 124
$ws.Range('K6').Value = 168
$ws.Range('K7').Value = 454
$ws = $wb.Worksheets.Item('Austin')
$ws.Range('K2').Value = 325
$ws.Range('K3').Value = 346
$ws.Range('K4').Value = 67
$ws.Range('K5').Value = 35
$ws.Range('K6').Value = 397
$ws.Range('K7').Value = 1170
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('K2').Value = 202
$ws.Range('K3').Value = 276
$ws.Range('K6').Value = 215
$ws.Range('K7').Value = 741
$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('K2').Value = 98
$ws.Range('K6').Value = 69
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('K2').Value = 170
$ws.Range('K3').Value = 192
$ws.Range('K7').Value = 589
$ws = $wb.Worksheets.Item('New City')
$ws.Range('K2').Value = 125
$ws.Range('K7').Value = 404
$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('K2').Value = 65
$ws.Range('K7').Value = 191
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('K2').Value = 194
$ws.Range('K3').Value = 235
$ws.Range('K7').Value = 671
$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('K2').Value = 18
$ws.Range('K6').Value = 52
$ws.Range('K7').Value = 99
$ws = $wb.Worksheets.Item('Loop')
$ws.Range('K6').Value = 180
$ws.Range('K7').Value = 342
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('K2').Value = 271
$ws.Range('K3').Value = 335
$ws.Range('K5').Value = 26
$ws.Range('K7').Value = 940
$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('K6').Value = 64
$ws.Range('K7').Value = 155
$ws = $wb.Worksheets.Item('River North')
$ws.Range('K6').Value = 131
$ws.Range('K7').Value = 241
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('K2').Value = 172
$ws.Range('K7').Value = 647
$ws = $wb.Worksheets.Item('Boystown')
$ws.Range('K2').Value = 3
$ws.Range('K6').Value = 19
$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('K6').Value = 46
$ws.Range('K7').Value = 95
$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('K2').Value = 58
$ws.Range('K6').Value = 72
$ws.Range('K7').Value = 201
$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('K6').Value = 70
$ws.Range('K7').Value = 199
$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('K6').Value = 47
$ws.Range('K7').Value = 176
$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('K3').Value = 87
$ws.Range('K7').Value = 187
$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('K2').Value = 133
$ws.Range('K3').Value = 128
$ws.Range('K6').Value = 114
$ws.Range('K7').Value = 399
$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('K2').Value = 31
$ws.Range('K4').Value = 15
$ws.Range('K7').Value = 118
$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('K5').Value = 4
$ws.Range('K7').Value = 230
$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('K6').Value = 29
$ws.Range('K7').Value = 94
$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('K3').Value = 43
$ws.Range('K7').Value = 231
$ws = $wb.Worksheets.Item('East Side')
$ws.Range('K2').Value = 27
$ws.Range('K7').Value = 83
$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('K2').Value = 61
$ws.Range('K7').Value = 176
$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('K2').Value = 22
$ws.Range('K6').Value = 42
$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('K6').Value = 56
$ws.Range('K7').Value = 151
$ws = $wb.Worksheets.Item('Montclare')
$ws.Range('K3').Value = 10
$ws.Range('K7').Value = 29
$ws = $wb.Worksheets.Item('United Center')
$ws.Range('K4').Value = 4
$ws.Range('K6').Value = 79
$ws.Range('K7').Value = 194
$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('K6').Value = 61
$ws.Range('K7').Value = 161
$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('K2').Value = 55
$ws.Range('K7').Value = 155
$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('K2').Value = 53
$ws.Range('K7').Value = 125
$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range('K2').Value = 16
$ws.Range('K7').Value = 63
$ws = $wb.Worksheets.Item('Sauganash,Forest Glen')
$ws.Range('K2').Value = 4
$ws.Range('K7').Value = 12
